$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.527.32'
$ws.Range('E2').Value = '  +4.73%  '
$ws.Range('D3').Value = '2.990.18'
$ws.Range('E3').Value = '  +4.80%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'506.27"
$ws.Range('E5').Value = '  +7.44%  '
$ws.Range('D6').Value = "'136.32"
$ws.Range('E6').Value = '  +8.97%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = "'0.432"
$ws.Range('E8').Value = '  +7.73%  '
$ws.Range('D9').Value = "'7.49"
$ws.Range('E9').Value = '  +13.54%  '
$ws.Range('E10').Value = '  +13.35%  '
$ws.Range('E11').Value = '  +7.53%  '
$ws.Range('E12').Value = '  +4.46%  '
$ws.Range('D13').Value = '3.498.99'
$ws.Range('E13').Value = '  +4.65%  '
$ws.Range('D14').Value = "'25.77"
$ws.Range('E14').Value = '  +12.79%  '
$ws.Range('D15').Value = "'0.0000154"
$ws.Range('E15').Value = '  +16.16%  '
$ws.Range('D16').Value = '56.531.46'
$ws.Range('E16').Value = '  +4.36%  '
$ws.Range('D17').Value = '2.985.14'
$ws.Range('E17').Value = '  +4.36%  '
$ws.Range('D18').Value = "'5.81"
$ws.Range('E18').Value = '  +9.99%  '
$ws.Range('D19').Value = "'12.46"
$ws.Range('E19').Value = '  +9.11%  '
$ws.Range('D20').Value = "'7.81"
$ws.Range('E20').Value = '  +11.59%  '
$ws.Range('D21').Value = "'326.56"
$ws.Range('E21').Value = '  +11.65%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = "'0.475"
$ws.Range('E23').Value = '  +7.63%  '
$ws.Range('D24').Value = "'62.30"
$ws.Range('E24').Value = '  +6.68%  '
$ws.Range('D25').Value = "'0.997"
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('E26').Value = '  +7.48%  '
$ws.Range('D27').Value = '0.0₃0905'
$ws.Range('E27').Value = '  +13.09%  '
$ws.Range('D28').Value = "'6.47"
$ws.Range('E28').Value = '  +4.07%  '
$ws.Range('D29').Value = "'6.95"
$ws.Range('E29').Value = '  +13.88%  '
$ws.Range('D30').Value = "'1.21"
$ws.Range('E30').Value = '  +8.86%  '
$ws.Range('E31').Value = '  +9.94%  '
$ws.Range('D32').Value = "'20.59"
$ws.Range('E32').Value = '  +9.13%  '
$ws.Range('D33').Value = "'157.27"
$ws.Range('E33').Value = '  +16.82%  '
$ws.Range('D34').Value = "'4.49"
$ws.Range('E34').Value = '  +7.42%  '
$ws.Range('D35').Value = "'5.58"
$ws.Range('E35').Value = '  +3.86%  '
$ws.Range('E36').Value = '  +5.72%  '
$ws.Range('D37').Value = "'0.0676"
$ws.Range('E37').Value = '  +11.15%  '
$ws.Range('D38').Value = "'23.41"
$ws.Range('E38').Value = '  +3.19%  '
$ws.Range('D39').Value = '3.022.23'
$ws.Range('E39').Value = '  +5.02%  '
$ws.Range('D40').Value = "'36.35"
$ws.Range('E40').Value = '  +3.99%  '
$ws.Range('D41').Value = "'1.00"
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('E42').Value = '  +8.23%  '
$ws.Range('D43').Value = '2.259.86'
$ws.Range('E43').Value = '  +11.32%  '
$ws.Range('D44').Value = "'1.41"
$ws.Range('E44').Value = '  +8.45%  '
$ws.Range('D45').Value = "'0.990"
$ws.Range('E45').Value = '  +4.20%  '
$ws.Range('D46').Value = "'3.59"
$ws.Range('E46').Value = '  +6.62%  '
$ws.Range('D47').Value = "'1.96"
$ws.Range('E47').Value = '  +23.21%  '
$ws.Range('E48').Value = '  +12.40%  '
$ws.Range('E49').Value = '  +9.71%  '
$ws.Range('D50').Value = "'19.15"
$ws.Range('E50').Value = '  +8.71%  '
$ws.Range('E51').Value = '  +11.43%  '
